$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data 20s")
$ws2 = $wb.Worksheets.Item("Data 60s")

# ---------------------------------------------------------------------
# Sheet "Data 20s": add Male/Female/Ratio summary formulas (rows 16-18)
# ---------------------------------------------------------------------
$ws1.Range("M16").Value = "Numer of M"
$ws1.Range("O16").Formula = "=SUMPRODUCT(LEN(D4:D63)-LEN(SUBSTITUTE(D4:D63,""M"","""")))"

$ws1.Range("M17").Value = "Number of F"
$ws1.Range("O17").Formula = "=SUMPRODUCT(LEN(D4:D63)-LEN(SUBSTITUTE(D4:D63,""F"","""")))"

$ws1.Range("M18").Value = "Ratio F/M"
$ws1.Range("O18").Formula = "=O17/O16"

# ---------------------------------------------------------------------
# Sheet "Data 20s": fill in the newly-recorded subjects (rows 46-54)
# ---------------------------------------------------------------------
$ws1.Range("A46").Value = "민정"
$ws1.Range("B46").Value = "김"
$ws1.Range("C46").Value = 1993
$ws1.Range("D46").Value = "F"
$ws1.Range("E46").Value = "No"

$ws1.Range("A47").Value = "민정"
$ws1.Range("B47").Value = "김"
$ws1.Range("C47").Value = 1993
$ws1.Range("D47").Value = "F"
$ws1.Range("E47").Value = "No"

$ws1.Range("A48").Value = "민정"
$ws1.Range("B48").Value = "김"
$ws1.Range("C48").Value = 1993
$ws1.Range("D48").Value = "F"
$ws1.Range("E48").Value = "No"

$ws1.Range("A49").Value = "은아"
$ws1.Range("B49").Value = "박"
$ws1.Range("C49").Value = 1993
$ws1.Range("D49").Value = "F"
$ws1.Range("E49").Value = "No"

$ws1.Range("A50").Value = "은아"
$ws1.Range("B50").Value = "박"
$ws1.Range("C50").Value = 1993
$ws1.Range("D50").Value = "F"
$ws1.Range("E50").Value = "No"

$ws1.Range("A51").Value = "은아"
$ws1.Range("B51").Value = "박"
$ws1.Range("C51").Value = 1993
$ws1.Range("D51").Value = "F"
$ws1.Range("E51").Value = "No"

$ws1.Range("A52").Value = "신희"
$ws1.Range("B52").Value = "박"
$ws1.Range("C52").Value = 1995
$ws1.Range("D52").Value = "F"
$ws1.Range("E52").Value = "No"

$ws1.Range("A53").Value = "신희"
$ws1.Range("B53").Value = "박"
$ws1.Range("C53").Value = 1995
$ws1.Range("D53").Value = "F"
$ws1.Range("E53").Value = "No"

$ws1.Range("A54").Value = "신희"
$ws1.Range("B54").Value = "박"
$ws1.Range("C54").Value = 1995
$ws1.Range("D54").Value = "F"
$ws1.Range("E54").Value = "No"

# Minor border tweak that trails along on rows 59-60, column D
# (bottom edge goes from medium back to thin)
$ws1.Range("D59").Borders.Item(9).Weight = 2
$ws1.Range("D60").Borders.Item(9).Weight = 2

# ---------------------------------------------------------------------
# Sheet "Data 60s": the Hand/Wrist IMU filename formulas referred to the
# wrong test ("..._20_...") - point them at the 60s test files instead.
# ---------------------------------------------------------------------
$ws2.Range("H4").Formula = "=CONCATENATE(""Hand_IMU_60_"",J3+1,"".txt"")"
$ws2.Range("I4").Formula = "=CONCATENATE(""Wrist_IMU_60_"",J3+1,"".txt"")"

$ws2.Range("H5").Formula = "=CONCATENATE(""Hand_IMU_60_"",J4+1,"".txt"")"
$ws2.Range("I5").Formula = "=CONCATENATE(""Wrist_IMU_60_"",J4+1,"".txt"")"

for ($r = 6; $r -le 22; $r++) {
    $ws2.Range("H$r").Formula = "=CONCATENATE(""Hand_IMU_60_"",J" + ($r - 1) + "+1,"".txt"")"
    $ws2.Range("I$r").Formula = "=CONCATENATE(""Wrist_IMU_60_"",J" + ($r - 1) + "+1,"".txt"")"
}

# ---------------------------------------------------------------------
# Sheet selections - "Data 20s" ends up the active tab/cell, "Data 60s"
# keeps its own last-used selection.
# ---------------------------------------------------------------------
$ws2.Range("I25").Select()
$ws1.Range("G14").Select()
